# Apply the "Upload new version with timestamp" changes to the daily
# missing-items sheet.
#
# Row 26 (TEMPO COOL): current balance, sale price and transactions count
# were recalculated.
# Row 27 (VISCERALGINE ...): order-limit value dropped from 1 to 0.
# Row 32 (سرنجات 5 سم): sale price and transactions count were recalculated.
# Row 37: total sale price recomputed.
# Row 38 (A38): footer timestamp updated to reflect the new export time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some cells hold numeric-looking text (e.g. "34.1900") in the
# source workbook even though their number format looks numeric. Setting
# .Value directly on those would make Excel coerce the text into a real
# number and lose the formatting/trailing zeros, so we temporarily force
# a text number format while assigning, then restore the original format.
function Set-TextValue($range, $value) {
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $origFormat
}

# Row 26 - TEMPO COOL: balance, sale price and transaction count updated
$ws.Range("H26").Value = "1:1"
Set-TextValue $ws.Range("P26") "34.1900"
$ws.Range("Q26").Value = "1:0"

# Row 27 - VISCERALGINE ...: order limit dropped from 1 to 0
Set-TextValue $ws.Range("L27") "0"

# Row 32 - سرنجات 5 سم: sale price and transaction count updated
Set-TextValue $ws.Range("P32") "6.0000"
$ws.Range("Q32").Value = "2:0"

# Row 37 - recomputed total of the sale price column
$ws.Range("P37").Value = 1077.915

# Row 38 - footer timestamp reflecting the new export time
$ws.Range("A38").Value = "Monday, 28 July, 2025 12:02 PM"
